$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$ws.Activate()

# Row 1 header: remove bold/alignment/wrap styling (back to default "Normal") and
# drop the explicit 30pt row height (back to default row height).
$ws.Range("A1:H1").Style = "Normal"
$ws.Rows.Item(1).AutoFit()

# Row 2 (LDVs) calibration values: 0.076 -> 0.075 across all vehicle-type columns.
$ws.Range("B2:H2").Value = 0.075

# Update the current selection to match the saved view.
$ws.Range("E15").Select()
